# Apply weekly update: insert a new data row at row 7 (pushing existing rows down)
# and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7; existing rows 7.. shift down to 8..
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new observation.
$ws.Range("A7").Value = 9
$ws.Range("B7").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C7").Value = "Metropolitana"
$ws.Range("D7").Value = 44552
$ws.Range("E7").Value = 13
$ws.Range("F7").Value = 100112029
$ws.Range("G7").Value = "Orégano"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 8
$ws.Range("K7").Value = 9000
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = 10000
$ws.Range("N7").Value = '$/docena de atados'
$ws.Range("O7").Value = "Región Metropolitana"
$ws.Range("P7").Value = 3333
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = "Hortaliza"
